# Update gh-pages output data (想去人数 / 最低票价 columns) to match newly
# generated values, per commit "Update gh-pages to output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# ---- Sheet "展览" (exhibitions) ----
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 80
$ws1.Range("F3").Value = 142
$ws1.Range("F4").Value = 2103
$ws1.Range("G4").Value = 55.2
$ws1.Range("F5").Value = 378
$ws1.Range("F6").Value = 662
$ws1.Range("F7").Value = 109
$ws1.Range("F8").Value = 2092
$ws1.Range("F9").Value = 10862
$ws1.Range("F15").Value = 9066
$ws1.Range("F17").Value = 737
$ws1.Range("F18").Value = 5323
$ws1.Range("F19").Value = 75
$ws1.Range("F20").Value = 3372

# ---- Sheet "演出" (performances) ----
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 559

# ---- Sheet "全部类型" (all types, combined) ----
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 80
$ws4.Range("F3").Value = 142
$ws4.Range("F4").Value = 2103
$ws4.Range("G4").Value = 55.2
$ws4.Range("F5").Value = 378
$ws4.Range("F6").Value = 662
$ws4.Range("F8").Value = 109
$ws4.Range("F9").Value = 2092
$ws4.Range("F10").Value = 559
$ws4.Range("F12").Value = 10862
$ws4.Range("F18").Value = 9066
$ws4.Range("F20").Value = 737
$ws4.Range("F21").Value = 5323
$ws4.Range("F22").Value = 75
$ws4.Range("F23").Value = 3372
